$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Food-Beverages" (sheet1.xml): append 4 new tickers (rows 30-33)
# ---------------------------------------------------------------------------
$wsFood = $wb.Worksheets.Item("Food-Beverages")

$wsFood.Cells.Item(30,1).Value = "x"
$wsFood.Cells.Item(30,2).Value = "Givaudan"
$wsFood.Cells.Item(30,3).Value = "GIVN SW"

$wsFood.Cells.Item(31,1).Value = "x"
$wsFood.Cells.Item(31,2).Value = "Tyson Foods"
$wsFood.Cells.Item(31,3).Value = "TSN"

$wsFood.Cells.Item(32,1).Value = "x"
$wsFood.Cells.Item(32,2).Value = "International Flavor"
$wsFood.Cells.Item(32,3).Value = "IFF"

$wsFood.Cells.Item(33,1).Value = "x"
$wsFood.Cells.Item(33,2).Value = "Hormel Foods"
$wsFood.Cells.Item(33,3).Value = "HRL"

[void]$wsFood.Activate()
[void]$wsFood.Range("D33").Select()

# ---------------------------------------------------------------------------
# Sheet "Apparel" (sheet2.xml): append 1 new ticker (row 13)
# ---------------------------------------------------------------------------
$wsApparel = $wb.Worksheets.Item("Apparel")

$wsApparel.Cells.Item(13,1).Value = "x"
$wsApparel.Cells.Item(13,2).Value = "Adidas"
$wsApparel.Cells.Item(13,3).Value = "ADS GR"

[void]$wsApparel.Activate()
[void]$wsApparel.Range("B14").Select()

# ---------------------------------------------------------------------------
# Sheet "Retail" (sheet3.xml): append 4 new tickers (rows 23-26)
# ---------------------------------------------------------------------------
$wsRetail = $wb.Worksheets.Item("Retail")

$wsRetail.Cells.Item(23,1).Value = "x"
$wsRetail.Cells.Item(23,2).Value = "Woolworths"
$wsRetail.Cells.Item(23,3).Value = "WOW AU"

$wsRetail.Cells.Item(24,1).Value = "x"
$wsRetail.Cells.Item(24,2).Value = "Loblaws"
$wsRetail.Cells.Item(24,3).Value = "L CN"

$wsRetail.Cells.Item(25,1).Value = "x"
$wsRetail.Cells.Item(25,2).Value = "Ahold"
$wsRetail.Cells.Item(25,3).Value = "AD NA"

$wsRetail.Cells.Item(26,1).Value = "x"
$wsRetail.Cells.Item(26,2).Value = "Ross Stores"
$wsRetail.Cells.Item(26,3).Value = "ROST"

[void]$wsRetail.Activate()
[void]$wsRetail.Range("B27").Select()

# ---------------------------------------------------------------------------
# Sheet "Restaurants" (sheet4.xml): no data changes, just re-select the
# header row as the last-used range.
# ---------------------------------------------------------------------------
$wsRest = $wb.Worksheets.Item("Restaurants")
[void]$wsRest.Activate()
[void]$wsRest.Range("B2:E2").Select()

# ---------------------------------------------------------------------------
# Sheet "Leisure" (sheet5.xml): rebuild to match the standard layout used by
# the other sheets (header row + "x" flag column + frozen panes) and add two
# new holdings (Las Vegas Sands, Galaxy Entertainment).
# ---------------------------------------------------------------------------
$wsLeisure = $wb.Worksheets.Item("Leisure")

# Shift the two existing data rows down one row (2->3, 3->4) and flag them.
$wsLeisure.Cells.Item(4,2).Value = $wsLeisure.Cells.Item(3,2).Value2
$wsLeisure.Cells.Item(4,3).Value = $wsLeisure.Cells.Item(3,3).Value2
$wsLeisure.Cells.Item(3,2).Value = $wsLeisure.Cells.Item(2,2).Value2
$wsLeisure.Cells.Item(3,3).Value = $wsLeisure.Cells.Item(2,3).Value2

$wsLeisure.Cells.Item(3,1).Value = "x"
$wsLeisure.Cells.Item(4,1).Value = "x"

# New header row, same as the other sheets.
$wsLeisure.Cells.Item(2,2).Value = "Name"
$wsLeisure.Cells.Item(2,3).Value = "Ticker"
$wsLeisure.Cells.Item(2,4).Value = "Price"
$wsLeisure.Cells.Item(2,5).Value = "MC"

# Two new holdings.
$wsLeisure.Cells.Item(5,1).Value = "x"
$wsLeisure.Cells.Item(5,2).Value = "Las Vegas Sands"
$wsLeisure.Cells.Item(5,3).Value = "LVS"

$wsLeisure.Cells.Item(6,1).Value = "x"
$wsLeisure.Cells.Item(6,2).Value = "Galaxy Entertainment"
$wsLeisure.Cells.Item(6,3).Value = "27 HK"

# Freeze the header (same split as every other sheet) and size the columns
# to fit their (now longer) contents, same as the other tabs' bestFit cols.
[void]$wsLeisure.Activate()
[void]$wsLeisure.Range("C3").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsLeisure.Columns.Item(1).ColumnWidth = 4.166666666666667
$wsLeisure.Columns.Item(2).ColumnWidth = 18.307291666666668
[void]$wsLeisure.Range("B7").Select()

# ---------------------------------------------------------------------------
# Sheet "NonDurable" (sheet6.xml): append 1 new ticker (row 14). This sheet
# becomes the active tab.
# ---------------------------------------------------------------------------
$wsNonDur = $wb.Worksheets.Item("NonDurable")

$wsNonDur.Cells.Item(14,1).Value = "x"
$wsNonDur.Cells.Item(14,2).Value = "Beiersdorf"
$wsNonDur.Cells.Item(14,3).Value = "BEI GR"

[void]$wsNonDur.Activate()
[void]$wsNonDur.Range("B15").Select()
